$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells C2:E2 / C4:E4 hold numbers stored as text ("runs","balls","fours").
# Force text-formatted assignment (matches the source data's
# numberStoredAsText convention) so the swapped values keep their string type.
$ws.Range("C2:E2").NumberFormat = "@"
$ws.Range("C4:E4").NumberFormat = "@"

# Row 2 (Karun Nair's first Kings XI Punjab innings row) takes the values
# that used to live in row 4 ...
$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "3"
$ws.Range("E2").Value = "0"

# ... and row 4 takes the values that used to live in row 2.
$ws.Range("C4").Value = "15"
$ws.Range("D4").Value = "8"
$ws.Range("E4").Value = "2"
